# Insert a new weekly data row at row 42 (pushes existing rows 42-123 down to 43-124)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(42).Insert()

$ws.Range("A42").Value = 6
$ws.Range("B42").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C42").Value = "Metropolitana"
$ws.Range("D42").Value = 44519
$ws.Range("E42").Value = 13
$ws.Range("F42").Value = 100112029
$ws.Range("G42").Value = "Orégano"
$ws.Range("H42").Value = "Sin especificar"
$ws.Range("I42").Value = "Primera"
$ws.Range("J42").Value = 33
$ws.Range("K42").Value = 8000
$ws.Range("L42").Value = 9000
$ws.Range("M42").Value = 8455
$ws.Range("N42").Value = "`$/docena de atados"
$ws.Range("O42").Value = "Región Metropolitana"
$ws.Range("P42").Value = 2818
$ws.Range("Q42").Value = 3
$ws.Range("R42").Value = "Hortaliza"
